$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format price cells whose new values look numeric so Excel keeps them
# as plain text, matching the source data (every Price-column cell is text).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = '61.892.81'
$ws.Range("E2").Value = '  -0.39%  '

$ws.Range("D3").Value = '3.414.42'
$ws.Range("E3").Value = '  -0.26%  '

$ws.Range("E4").Value = '  +0.38%  '

$ws.Range("D5").Value = '408.89'
$ws.Range("E5").Value = '  +0.21%  '

$ws.Range("D6").Value = '128.66'
$ws.Range("E6").Value = '  -3.46%  '

$ws.Range("D7").Value = '0.631'
$ws.Range("E7").Value = '  +6.38%  '

$ws.Range("D9").Value = '0.743'
$ws.Range("E9").Value = '  +9.80%  '

$ws.Range("D10").Value = '0.142'
$ws.Range("E10").Value = '  +16.04%  '

$ws.Range("D11").Value = '42.59'
$ws.Range("E11").Value = '  +0.59%  '

$ws.Range("D12").Value = '0.0000215'
$ws.Range("E12").Value = '  +62.56%  '

$ws.Range("E13").Value = '  -0.33%  '

$ws.Range("D14").Value = '3.969.39'
$ws.Range("E14").Value = '  +0.02%  '

$ws.Range("D15").Value = '8.92'
$ws.Range("E15").Value = '  +5.52%  '

$ws.Range("D16").Value = '20.94'
$ws.Range("E16").Value = '  +4.73%  '

$ws.Range("D17").Value = '3.364.56'
$ws.Range("E17").Value = '  -1.84%  '

$ws.Range("D18").Value = '12.17'
$ws.Range("E18").Value = '  +10.30%  '

$ws.Range("D19").Value = '1.06'
$ws.Range("E19").Value = '  +3.76%  '

$ws.Range("D20").Value = '61.957.86'
$ws.Range("E20").Value = '  -0.26%  '

$ws.Range("D21").Value = '400.34'
$ws.Range("E21").Value = '  +26.87%  '

$ws.Range("D22").Value = '89.46'
$ws.Range("E22").Value = '  +4.87%  '

$ws.Range("D23").Value = '3.18'
$ws.Range("E23").Value = '  -1.05%  '

$ws.Range("D24").Value = '13.12'
$ws.Range("E24").Value = '  +2.46%  '

$ws.Range("D25").Value = '3.24'
$ws.Range("E25").Value = '  +3.39%  '

$ws.Range("D26").Value = '32.57'
$ws.Range("E26").Value = '  +9.40%  '

$ws.Range("D27").Value = '8.59'
$ws.Range("E27").Value = '  +4.06%  '

$ws.Range("E28").Value = '  +0.52%  '

$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '2.74'
$ws.Range("E29").Value = '  +0.66%  '

$ws.Range("B30").Value = 'RenderToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D30").Value = '7.58'
$ws.Range("E30").Value = '  -1.88%  '

$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").Value = '0.120'
$ws.Range("E31").Value = '  +2.86%  '

$ws.Range("D32").Value = '0.171'
$ws.Range("E32").Value = '  -1.80%  '

$ws.Range("D33").Value = '11.85'
$ws.Range("E33").Value = '  +3.69%  '

$ws.Range("D34").Value = '43.12'
$ws.Range("E34").Value = '  +0.82%  '

$ws.Range("E35").Value = '  +0.70%  '

$ws.Range("D36").Value = '0.0496'
$ws.Range("E36").Value = '  +2.06%  '

$ws.Range("D37").Value = '54.04'
$ws.Range("E37").Value = '  +3.61%  '

$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  +0.11%  '

$ws.Range("D39").Value = '3.37'
$ws.Range("E39").Value = '  -1.90%  '

$ws.Range("D40").Value = '0.133'
$ws.Range("E40").Value = '  +7.14%  '

$ws.Range("D41").Value = '2.90'
$ws.Range("E41").Value = '  -2.96%  '

$ws.Range("B42").Value = 'TheGraph'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D42").Value = '0.310'
$ws.Range("E42").Value = '  +6.33%  '

$ws.Range("B43").Value = 'Monero'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D43").Value = '141.02'
$ws.Range("E43").Value = '  +2.13%  '

$ws.Range("D44").Value = '1.97'
$ws.Range("E44").Value = '  -1.81%  '

$ws.Range("D45").Value = '4.06'
$ws.Range("E45").Value = '  +1.44%  '

$ws.Range("D46").Value = '2.41'
$ws.Range("E46").Value = '  +8.90%  '

$ws.Range("D47").Value = '16.59'
$ws.Range("E47").Value = '  -1.67%  '

$ws.Range("D48").Value = '21.94'
$ws.Range("E48").Value = '  +1.91%  '

$ws.Range("D49").Value = '2.123.55'
$ws.Range("E49").Value = '  -0.42%  '

$ws.Range("D50").Value = '2.38'
$ws.Range("E50").Value = '  +3.73%  '

$ws.Range("E51").Value = '  +17.24%  '
